$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Region Dict")

# Rename the two "Region Norden N" entries to "Region Nord N"
$ws.Range("A3").Value = "Region Nord 1"
$ws.Range("A4").Value = "Region Nord 2"

# Move the active selection to A4 (matches saved selection in the file)
$ws.Activate()
$ws.Range("A4").Select()
